$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price (D) and Volume(1h) (E) columns for updated cells
# so Excel does not auto-convert numeric-looking strings into numbers,
# matching the original inline-string cell type used in the workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.091.43'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.302.30'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.18%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.93'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.88'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.508'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.34%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.81'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0794'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.36'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.42%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.29'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +14.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.78'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.656.95'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.330.01'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.809'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.044.56'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.76'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0902'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.07'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.66'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.80'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.65%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.05'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +6.80%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.48'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.63%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.17'

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.45'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.02'

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.10'

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.55%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.12%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +8.04%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.70%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0699'

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.84'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.88%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.28%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.93%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.39'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.981.79'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.34%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.96'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.52'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.95%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.535.19'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.30'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.21%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.50%  '
